$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2018-10", 102.1, 100.6),
    @("2018-11", 101.9, 99.9),
    @("2018-12", 101.2, 99.9),
    @("2018-01", 101.447, 104.2028),
    @("2018-02", 101.3, 102.7),
    @("2018-03", 101.2, 101.2),
    @("2018-04", 101.3, 100.5),
    @("2018-05", 101.9, 100.6),
    @("2018-06", 102.4, 101.1),
    @("2018-07", 102.6, 101),
    @("2018-08", 102.6, 100.9),
    @("2018-09", 102.1, 100.8),
    @("2019-10", 97.8, 100.1),
    @("2019-11", 97.8, 100),
    @("2019-12", 98.1, 99.9),
    @("2019-01", 100.7, 99.6),
    @("2019-02", 100.2, 99.4),
    @("2019-03", 100.1, 99.1),
    @("2019-04", 100.4, 99.4),
    @("2019-05", 100, 99.6),
    @("2019-06", 99.5, 99.9),
    @("2019-07", 99.2, 100),
    @("2019-08", 98.6, 100.1),
    @("2019-09", 98.3, 100.2),
    @("2020-10", 97.5, 97.5),
    @("2020-11", 98.2, 97.8),
    @("2020-12", 99.4, 98.2),
    @("2020-01", 98.3, 100),
    @("2020-02", 98.5, 100.1),
    @("2020-03", 98.1, 99.8),
    @("2020-04", 97.2, 99.5),
    @("2020-05", 96.9, 99.5),
    @("2020-06", 97.2, 99.2),
    @("2020-07", 97.1, 98.9),
    @("2020-08", 97.2, 98.3),
    @("2020-09", 97.3, 97.9),
    @("2021-10", 106.5, 103.8),
    @("2021-11", 106.9, 104.7),
    @("2021-12", 105.1, 104.7),
    @("2021-01", 99.3, 98.2),
    @("2021-02", 99.8, 98.5),
    @("2021-03", 101.6, 100.1),
    @("2021-04", 103.1, 100.9),
    @("2021-05", 103.7, 101.7),
    @("2021-06", 104, 101.6),
    @("2021-07", 104.3, 102.4),
    @("2021-08", 104.6, 102.9),
    @("2021-09", 105.1, 102.9),
    @("2022-10", 97.7, 101.7),
    @("2022-11", 96.7, 101.4),
    @("2022-12", 96.9, 101.1),
    @("2022-01", 104.4, 104.1),
    @("2022-02", 104.1, 103.9),
    @("2022-03", 103.2, 103),
    @("2022-04", 102.8, 102.6),
    @("2022-05", 102.7, 101.9),
    @("2022-06", 102.4, 102.2),
    @("2022-07", 101.5, 102.1),
    @("2022-08", 100.2, 101.4),
    @("2022-09", 99.5, 101.8),
    @("2023-01", 96.7, 100.7),
    @("2023-02", 96.5, 99.9),
    @("2023-03", 96, 99.8),
    @("2023-04", 95.1, 99.4),
    @("2023-05", 94.8, 99.1),
    @("2023-06", 94.5, 98.7),
    @("2023-07", 94.9, 98.5)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Make sure the date column keeps the same (bold / bordered / centered)
# style used by the existing rows, including the newly appended ones.
$lastRow = $r - 1
$ws.Range("A2").Copy()
$ws.Range("A2:A" + $lastRow).PasteSpecial(-4122)

Write-Output "updated $($rows.Length) rows"
